$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the time_taken column, matching style of existing headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate the time_taken values for each data row as text strings
$times = @(
    "2021-10-05 13:41:06.542439",
    "2021-10-05 13:41:06.542451",
    "2021-10-05 13:41:06.542455",
    "2021-10-05 13:41:06.542458",
    "2021-10-05 13:41:06.542461",
    "2021-10-05 13:41:06.542465",
    "2021-10-05 13:41:06.542468",
    "2021-10-05 13:41:06.542471",
    "2021-10-05 13:41:06.542474",
    "2021-10-05 13:41:06.542477",
    "2021-10-05 13:41:06.542480"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
